$wb = $excel.ActiveWorkbook

# --- Sheet "Overview": Correspond Handoff Datetime style column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-12-16 09:56:15"
$wsOverview.Range("G3").Value = "2016-12-16 09:56:15"

# --- Sheet "zh-cn" ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "mt"
$wsZh.Range("E3").Value = "mt"
$wsZh.Range("H2").Value = "2016-12-16 09:56:01"
$wsZh.Range("H3").Value = "2016-12-16 09:56:01"
$wsZh.Range("L2").Value = "2016-12-16 09:56:55"
$wsZh.Range("L3").Value = "2016-12-16 09:56:55"

# --- Sheet "de-de" ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "mt"
$wsDe.Range("E3").Value = "mt"
$wsDe.Range("H2").Value = "2016-12-16 09:56:15"
$wsDe.Range("H3").Value = "2016-12-16 09:56:15"
$wsDe.Range("L2").Value = "2016-12-16 09:57:14"
$wsDe.Range("L3").Value = "2016-12-16 09:57:14"
